# Auto-generated Excel COM-interop edit script
# Applies updated crypto price/volume data scraped on Tue Aug  1 04:36:25 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.868.67"
$ws.Range("E2").Value = "  -1.91%  "

$ws.Range("D3").Value = "1.823.68"
$ws.Range("E3").Value = "  -2.33%  "

$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'238.88"
$ws.Range("E5").Value = "  -1.87%  "

$ws.Range("D6").Value = "'0.6893"
$ws.Range("E6").Value = "  -2.06%  "

$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").Value = "'0.07605"
$ws.Range("E8").Value = "  -3.02%  "

$ws.Range("D9").Value = "'0.3012"
$ws.Range("E9").Value = "  -3.98%  "

$ws.Range("D10").Value = "'23.36"
$ws.Range("E10").Value = "  -4.43%  "

$ws.Range("E11").Value = "  -3.70%  "

$ws.Range("D12").Value = "1.826.69"
$ws.Range("E12").Value = "  -2.96%  "

$ws.Range("D13").Value = "'5.038"
$ws.Range("E13").Value = "  -2.90%  "

$ws.Range("D14").Value = "'89.91"
$ws.Range("E14").Value = "  -3.59%  "

$ws.Range("D15").Value = "'0.6706"
$ws.Range("E15").Value = "  -4.18%  "

$ws.Range("D16").Value = "'6.410"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").Value = "'0.000008265"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").Value = "28.852.04"
$ws.Range("E18").Value = "  -2.18%  "

$ws.Range("D19").Value = "'242.61"
$ws.Range("E19").Value = "  -4.44%  "

$ws.Range("D20").Value = "2.080.56"
$ws.Range("E20").Value = "  -2.88%  "

$ws.Range("D21").Value = "'12.59"
$ws.Range("E21").Value = "  -4.04%  "

$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = "  -0.19%  "

$ws.Range("D23").Value = "'7.365"
$ws.Range("E23").Value = "  -2.98%  "

$ws.Range("D24").Value = "'0.9999"
$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").Value = "'0.1471"
$ws.Range("E25").Value = "  -5.26%  "

$ws.Range("D26").Value = "'160.24"
$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("D27").Value = "'8.694"
$ws.Range("E27").Value = "  -3.67%  "

$ws.Range("D28").Value = "'18.13"
$ws.Range("E28").Value = "  -3.32%  "

$ws.Range("D29").Value = "'1.526"
$ws.Range("E29").Value = "  +1.71%  "

$ws.Range("E30").Value = "  -3.24%  "

$ws.Range("D31").Value = "'4.127"
$ws.Range("E31").Value = "  -2.93%  "

$ws.Range("D32").Value = "'1.190"
$ws.Range("E32").Value = "  -0.96%  "

$ws.Range("D33").Value = "'0.05092"
$ws.Range("E33").Value = "  -3.91%  "

$ws.Range("D34").Value = "'0.7500"
$ws.Range("E34").Value = "  +0.76%  "

$ws.Range("D35").Value = "'1.809"
$ws.Range("E35").Value = "  -3.99%  "

$ws.Range("D36").Value = "'1.138"
$ws.Range("E36").Value = "  -2.21%  "

$ws.Range("E37").Value = "  -1.41%  "

$ws.Range("D38").Value = "'0.01826"
$ws.Range("E38").Value = "  -2.50%  "

$ws.Range("D39").Value = "1.200.69"
$ws.Range("E39").Value = "  -4.55%  "

$ws.Range("D40").Value = "'2.673"
$ws.Range("E40").Value = "  -2.54%  "

$ws.Range("D41").Value = "'0.9130"
$ws.Range("E41").Value = "  +1.71%  "

$ws.Range("D42").Value = "'107.97"
$ws.Range("E42").Value = "  -0.61%  "

$ws.Range("D43").Value = "'0.9990"
$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("D44").Value = "1.980.67"
$ws.Range("E44").Value = "  -2.85%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.00000000122"
$ws.Range("E45").Value = "  -6.31%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.5150"
$ws.Range("E46").Value = "  -0.77%  "

$ws.Range("D47").Value = "'9.415"
$ws.Range("E47").Value = "  -0.63%  "

$ws.Range("D48").Value = "'5.225"
$ws.Range("E48").Value = "  -12.05%  "

$ws.Range("D49").Value = "'1.723"
$ws.Range("E49").Value = "  -3.80%  "

$ws.Range("D50").Value = "'62.17"
$ws.Range("E50").Value = "  -12.67%  "

$ws.Range("D51").Value = "'0.4186"
$ws.Range("E51").Value = "  -2.66%  "
